# Applies the LOQ4243 syllabus content update described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação: date bumped from 2018 to 2021
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# Objetivos: (PT) / Objectives: (EN)
$ws.Range("B10").Value = "Este curso visa apresentar os principais conceitos relacionados ao uso da eletrônica, com foco na identificação de componentes e funcionalidades, para o uso em projetos relacionados ao curso de graduação em engenharia de produção."
$ws.Range("C10").Value = "Este curso visa apresentar os principais conceitos relacionados ao uso da eletrônica, com foco na identificação de componentes e funcionalidades, para o uso em projetos relacionados ao curso de graduação em engenharia de produção."

$ws.Range("B11").Value = "This course aims to present the main concepts related to the use of electronics, focusing on the identification of components and functionalities, for use in projects related to the undergraduate course in production engineering"
$ws.Range("C11").Value = "This course aims to present the main concepts related to the use of electronics, focusing on the identification of components and functionalities, for use in projects related to the undergraduate course in production engineering"

# Programa resumido: (PT) / Short syllabus: (EN)
$ws.Range("B14").Value = "Fundamentos da Eletrônica; Sensores e atuadores; Microcontroladores;Sistemas mecatrônicos;Experimentos práticos."
$ws.Range("C14").Value = "Fundamentos da Eletrônica; Sensores e atuadores; Microcontroladores;Sistemas mecatrônicos;Experimentos práticos."

$ws.Range("B15").Value = "Fundamentals of Electronics; Sensors and actuators; Microcontrollers; Mechatronic systems; Practical experiments."
$ws.Range("C15").Value = "Fundamentals of Electronics; Sensors and actuators; Microcontrollers; Mechatronic systems; Practical experiments."

# Programa: (PT) / Syllabus: (EN)
$ws.Range("B16").Value = "Fundamentos da Eletrônica: digital e analógica; Uso de sensores e atuadores em processos produtivos; Programação básica de microcontroladores, com foco em Arduino ou similar;Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;Experimentos práticos."
$ws.Range("C16").Value = "Fundamentos da Eletrônica: digital e analógica; Uso de sensores e atuadores em processos produtivos; Programação básica de microcontroladores, com foco em Arduino ou similar;Projetos de uso em sistemas mecatrônicos aplicados a engenharia de produção;Experimentos práticos."

$ws.Range("B17").Value = "Fundamentals of Electronics: digital and analog; Use of sensors and actuators in production processes; Basic programming of microcontrollers, focusing on Arduino or similar; Projects for use in mechatronic systems applied to production engineering; Practical experiments"
$ws.Range("C17").Value = "Fundamentals of Electronics: digital and analog; Use of sensors and actuators in production processes; Basic programming of microcontrollers, focusing on Arduino or similar; Projects for use in mechatronic systems applied to production engineering; Practical experiments"

# Método:
$ws.Range("B19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."
$ws.Range("C19").Value = "Esta disciplina deverá conter no mínimo duas avaliações denominadas A1 e A2. As avalições poderão ser: escritas, práticas, seminários, trabalhos de campo, projetos, ou outra forma de avaliação definida pelo professor."

# Critério:
$ws.Range("B20").Value = "Média ponderada das avaliações (M)."
$ws.Range("C20").Value = "Média ponderada das avaliações (M)."

# Norma de recuperação:
$ws.Range("B21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"
$ws.Range("C21").Value = "A recuperação será composta por uma única prova (RC) englobando toda a matéria ministrada ao longo do semestre. A média final, para os alunos em recuperação, será calculada com base na relação: MF=(M+RC)/2"

# Bibliografia:
$ws.Range("B22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.Bibliografia complementar será indicada ao longo do curso."
$ws.Range("C22").Value = "Tutoriais de Arduino disponibilizados pelo fabricante (arduino.cc) BOYLESTAD, Robert L.; NASHELSKY, Louis. Dispositivos Eletrônicos e Teoria de Circuitos. 8ª ed. São Paulo: Pearson. 696 p. THOMAZINI, Daniel; ALBUQUERQUE, Pedro U.B. Sensores Industriais – Fundamentos e Aplicações. 8ª ed. São Paulo: Érica, 2011. 224 p.Bibliografia complementar será indicada ao longo do curso."
